$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 13.13570827151799
$ws.Cells.Item(2, 4).Value = 4.909754846414156
$ws.Cells.Item(2, 5).Value = 13.38319263054962
$ws.Cells.Item(2, 6).Value = 66.36719626318511
$ws.Cells.Item(2, 7).Value = 91.92550057054183
$ws.Cells.Item(2, 8).Value = 29.17965981741694
$ws.Cells.Item(2, 10).Value = 12.35636330501284
$ws.Cells.Item(2, 12).Value = 9.292654946415571
$ws.Cells.Item(2, 14).Value = 18.99769123883691
$ws.Cells.Item(3, 3).Value = 13.13568159928572
$ws.Cells.Item(3, 4).Value = 4.839045006994266
$ws.Cells.Item(3, 5).Value = 13.41966388238334
$ws.Cells.Item(3, 6).Value = 65.86153012665315
$ws.Cells.Item(3, 7).Value = 90.82479694943839
$ws.Cells.Item(3, 8).Value = 29.05688509001541
$ws.Cells.Item(3, 10).Value = 12.38562875396479
$ws.Cells.Item(3, 12).Value = 9.317239672224913
$ws.Cells.Item(3, 14).Value = 18.40031508502705
$ws.Cells.Item(4, 3).Value = 13.13886024888667
$ws.Cells.Item(4, 4).Value = 4.794526842848682
$ws.Cells.Item(4, 5).Value = 13.44419776024397
$ws.Cells.Item(4, 6).Value = 65.57117214644053
$ws.Cells.Item(4, 7).Value = 90.17573754448048
$ws.Cells.Item(4, 8).Value = 28.9900010595173
$ws.Cells.Item(4, 10).Value = 12.40586554833572
$ws.Cells.Item(4, 12).Value = 9.333154588726861
$ws.Cells.Item(4, 14).Value = 18.02485520896364
$ws.Cells.Item(5, 3).Value = 13.14095568644013
$ws.Cells.Item(5, 4).Value = 4.776113773612114
$ws.Cells.Item(5, 5).Value = 13.45473318117156
$ws.Cells.Item(5, 6).Value = 65.45798041088217
$ws.Cells.Item(5, 7).Value = 89.91822218555495
$ws.Cells.Item(5, 8).Value = 28.96488851197993
$ws.Cells.Item(5, 10).Value = 12.41468085812786
$ws.Cells.Item(5, 12).Value = 9.339847152604325
$ws.Cells.Item(5, 14).Value = 17.86990355188116
$ws.Cells.Item(6, 3).Value = 13.14135185321532
$ws.Cells.Item(6, 4).Value = 4.77304006219359
$ws.Cells.Item(6, 5).Value = 13.45651503773418
$ws.Cells.Item(6, 6).Value = 65.43949682095432
$ws.Cells.Item(6, 7).Value = 89.87589053710894
$ws.Cells.Item(6, 8).Value = 28.96084815066821
$ws.Cells.Item(6, 10).Value = 12.41617893244457
$ws.Cells.Item(6, 12).Value = 9.34097098460412
$ws.Cells.Item(6, 14).Value = 17.84406337566596
$ws.Cells.Item(7, 3).Value = 13.13888527368356
$ws.Cells.Item(7, 4).Value = 4.794279609809277
$ws.Cells.Item(7, 5).Value = 13.4443376683211
$ws.Cells.Item(7, 6).Value = 65.56962473913063
$ws.Cells.Item(7, 7).Value = 90.17223601880124
$ws.Cells.Item(7, 8).Value = 28.98965369949104
$ws.Cells.Item(7, 10).Value = 12.40598213431083
$ws.Cells.Item(7, 12).Value = 9.333244007138218
$ws.Cells.Item(7, 14).Value = 18.02277304766463
$ws.Cells.Item(8, 3).Value = 13.13503406620736
$ws.Cells.Item(8, 4).Value = 4.885606090086736
$ws.Cells.Item(8, 5).Value = 13.39532334705613
$ws.Cells.Item(8, 6).Value = 66.18871196416556
$ws.Cells.Item(8, 7).Value = 91.5405745487906
$ws.Cells.Item(8, 8).Value = 29.13556723852023
$ws.Cells.Item(8, 10).Value = 12.36598242573276
$ws.Cells.Item(8, 12).Value = 9.300962192144336
$ws.Cells.Item(8, 14).Value = 18.79364780656866
$ws.Cells.Item(9, 3).Value = 13.15296658096233
$ws.Cells.Item(9, 4).Value = 5.055722423427215
$ws.Cells.Item(9, 5).Value = 13.3162153413174
$ws.Cells.Item(9, 6).Value = 67.5588022130955
$ws.Cells.Item(9, 7).Value = 94.42466347369628
$ws.Cells.Item(9, 8).Value = 29.48876504272338
$ws.Cells.Item(9, 10).Value = 12.30560763148381
$ws.Cells.Item(9, 12).Value = 9.244118589551229
$ws.Cells.Item(9, 14).Value = 20.2273683202997
$ws.Cells.Item(10, 3).Value = 13.18184017742235
$ws.Cells.Item(10, 4).Value = 5.174948579169153
$ws.Cells.Item(10, 5).Value = 13.26850148064461
$ws.Cells.Item(10, 6).Value = 68.65511171832949
$ws.Cells.Item(10, 7).Value = 96.64949959343994
$ws.Cells.Item(10, 8).Value = 29.78845194623611
$ws.Cells.Item(10, 10).Value = 12.27236206448919
$ws.Cells.Item(10, 12).Value = 9.206233211141303
$ws.Cells.Item(10, 14).Value = 21.22223697909767
$ws.Cells.Item(11, 3).Value = 13.19841109087211
$ws.Cells.Item(11, 4).Value = 5.227881757013342
$ws.Cells.Item(11, 5).Value = 13.24906354478717
$ws.Cells.Item(11, 6).Value = 69.17200953777981
$ws.Cells.Item(11, 7).Value = 97.68094682909805
$ws.Cells.Item(11, 8).Value = 29.93331183983019
$ws.Cells.Item(11, 10).Value = 12.25967186736832
$ws.Cells.Item(11, 12).Value = 9.189826973245108
$ws.Cells.Item(11, 14).Value = 21.66018057919901
$ws.Cells.Item(12, 3).Value = 13.20518160696949
$ws.Cells.Item(12, 4).Value = 5.247734323332162
$ws.Cells.Item(12, 5).Value = 13.2420296947837
$ws.Cells.Item(12, 6).Value = 69.37023878154591
$ws.Cells.Item(12, 7).Value = 98.07399808551676
$ws.Cells.Item(12, 8).Value = 29.98937197065196
$ws.Cells.Item(12, 10).Value = 12.25521805012015
$ws.Cells.Item(12, 12).Value = 9.183732382461947
$ws.Cells.Item(12, 14).Value = 21.82377585681918
$ws.Cells.Item(13, 3).Value = 13.20370140180643
$ws.Cells.Item(13, 4).Value = 5.243467331644287
$ws.Cells.Item(13, 5).Value = 13.24353000889991
$ws.Cells.Item(13, 6).Value = 69.32743778010743
$ws.Cells.Item(13, 7).Value = 97.98924292881092
$ws.Cells.Item(13, 8).Value = 29.97724520575006
$ws.Cells.Item(13, 10).Value = 12.25616159001053
$ws.Cells.Item(13, 12).Value = 9.185039724831014
$ws.Cells.Item(13, 14).Value = 21.78864458690801
$ws.Cells.Item(14, 3).Value = 13.19895816903234
$ws.Cells.Item(14, 4).Value = 5.22951891471879
$ws.Cells.Item(14, 5).Value = 13.24847830981323
$ws.Cells.Item(14, 6).Value = 69.18826877130317
$ws.Cells.Item(14, 7).Value = 97.71323582550512
$ws.Cells.Item(14, 8).Value = 29.93789992741326
$ws.Cells.Item(14, 10).Value = 12.25929839075908
$ws.Cells.Item(14, 12).Value = 9.189323206371105
$ws.Cells.Item(14, 14).Value = 21.67368539489661
$ws.Cells.Item(15, 3).Value = 13.19611736010382
$ws.Cells.Item(15, 4).Value = 5.22094995923114
$ws.Cells.Item(15, 5).Value = 13.25155188081935
$ws.Cells.Item(15, 6).Value = 69.10334443098414
$ws.Cells.Item(15, 7).Value = 97.54448499807063
$ws.Cells.Item(15, 8).Value = 29.91395602455493
$ws.Cells.Item(15, 10).Value = 12.26126562388254
$ws.Cells.Item(15, 12).Value = 9.191962313089309
$ws.Cells.Item(15, 14).Value = 21.60297336126124
$ws.Cells.Item(16, 3).Value = 13.18082649889512
$ws.Cells.Item(16, 4).Value = 5.171462639368021
$ws.Cells.Item(16, 5).Value = 13.26981749091202
$ws.Cells.Item(16, 6).Value = 68.62168666811473
$ws.Cells.Item(16, 7).Value = 96.58245435307026
$ws.Cells.Item(16, 8).Value = 29.7791548550043
$ws.Cells.Item(16, 10).Value = 12.27324051585281
$ws.Cells.Item(16, 12).Value = 9.207321976127862
$ws.Cells.Item(16, 14).Value = 21.19330956972086
$ws.Cells.Item(17, 3).Value = 13.172327230208
$ws.Cells.Item(17, 4).Value = 5.140766536579849
$ws.Cells.Item(17, 5).Value = 13.28160415744345
$ws.Cells.Item(17, 6).Value = 68.33077587364211
$ws.Cells.Item(17, 7).Value = 95.99701092761156
$ws.Cells.Item(17, 8).Value = 29.69862881748112
$ws.Cells.Item(17, 10).Value = 12.2812112315589
$ws.Cells.Item(17, 12).Value = 9.216956015557567
$ws.Cells.Item(17, 14).Value = 20.93814219015166
$ws.Cells.Item(18, 3).Value = 13.16776204962476
$ws.Cells.Item(18, 4).Value = 5.122988527991819
$ws.Cells.Item(18, 5).Value = 13.28859689783839
$ws.Cells.Item(18, 6).Value = 68.16517070828408
$ws.Cells.Item(18, 7).Value = 95.66212064698554
$ws.Cells.Item(18, 8).Value = 29.65311679593657
$ws.Cells.Item(18, 10).Value = 12.28602471241931
$ws.Cells.Item(18, 12).Value = 9.2225752549075
$ws.Cells.Item(18, 14).Value = 20.79000725568362
$ws.Cells.Item(19, 3).Value = 13.1662718497144
$ws.Cells.Item(19, 4).Value = 5.116948320296058
$ws.Cells.Item(19, 5).Value = 13.29100113954163
$ws.Cells.Item(19, 6).Value = 68.10939855761879
$ws.Cells.Item(19, 7).Value = 95.54905851210921
$ws.Cells.Item(19, 8).Value = 29.63784599671759
$ws.Cells.Item(19, 10).Value = 12.28769374024919
$ws.Cells.Item(19, 12).Value = 9.22449125796023
$ws.Cells.Item(19, 14).Value = 20.73962067985785
$ws.Cells.Item(20, 3).Value = 13.17319851278455
$ws.Cells.Item(20, 4).Value = 5.144046893981914
$ws.Cells.Item(20, 5).Value = 13.28032735859444
$ws.Cells.Item(20, 6).Value = 68.36156678765025
$ws.Cells.Item(20, 7).Value = 96.05914412413355
$ws.Cells.Item(20, 8).Value = 29.70711784796469
$ws.Cells.Item(20, 10).Value = 12.28033902842992
$ws.Cells.Item(20, 12).Value = 9.215922390486043
$ws.Cells.Item(20, 14).Value = 20.96544799483449
$ws.Cells.Item(21, 3).Value = 13.20033791579103
$ws.Cells.Item(21, 4).Value = 5.233621151472451
$ws.Cells.Item(21, 5).Value = 13.24701599585166
$ws.Cells.Item(21, 6).Value = 69.22907947623979
$ws.Cells.Item(21, 7).Value = 97.7942414559449
$ws.Cells.Item(21, 8).Value = 29.9494240762066
$ws.Cells.Item(21, 10).Value = 12.25836747766063
$ws.Cells.Item(21, 12).Value = 9.188061846696909
$ws.Cells.Item(21, 14).Value = 21.70751365554064
$ws.Cells.Item(22, 3).Value = 13.22096340698366
$ws.Cells.Item(22, 4).Value = 5.291041576425668
$ws.Cells.Item(22, 5).Value = 13.22715077172102
$ws.Cells.Item(22, 6).Value = 69.81050981515317
$ws.Cells.Item(22, 7).Value = 98.94244608913245
$ws.Cells.Item(22, 8).Value = 30.11479618349554
$ws.Cells.Item(22, 10).Value = 12.24605850294761
$ws.Cells.Item(22, 12).Value = 9.170541266988861
$ws.Cells.Item(22, 14).Value = 22.17935961385674
$ws.Cells.Item(23, 3).Value = 13.20969057099
$ws.Cells.Item(23, 4).Value = 5.260499294461462
$ws.Cells.Item(23, 5).Value = 13.2375785648729
$ws.Cells.Item(23, 6).Value = 69.49890823918572
$ws.Cells.Item(23, 7).Value = 98.32843190591954
$ws.Cells.Item(23, 8).Value = 30.02590027946225
$ws.Cells.Item(23, 10).Value = 12.25243981281001
$ws.Cells.Item(23, 12).Value = 9.179829709834952
$ws.Cells.Item(23, 14).Value = 21.92877110911181
$ws.Cells.Item(24, 3).Value = 13.17280360569008
$ws.Cells.Item(24, 4).Value = 5.142564250242428
$ws.Cells.Item(24, 5).Value = 13.28090392528953
$ws.Cells.Item(24, 6).Value = 68.34764108755151
$ws.Cells.Item(24, 7).Value = 96.03104842109684
$ws.Cells.Item(24, 8).Value = 29.70327751531613
$ws.Cells.Item(24, 10).Value = 12.28073263205317
$ws.Cells.Item(24, 12).Value = 9.21638944156266
$ws.Cells.Item(24, 14).Value = 20.95310750188673
$ws.Cells.Item(25, 3).Value = 13.14536869356388
$ws.Cells.Item(25, 4).Value = 5.010692816027857
$ws.Cells.Item(25, 5).Value = 13.3357914050063
$ws.Cells.Item(25, 6).Value = 67.17197704623268
$ws.Cells.Item(25, 7).Value = 93.62461101388271
$ws.Cells.Item(25, 8).Value = 29.38609650979222
$ws.Cells.Item(25, 10).Value = 12.31999639875615
$ws.Cells.Item(25, 12).Value = 9.258810910984556
$ws.Cells.Item(25, 14).Value = 19.84905939529497
